# Add 4 new RISPERDAL (Risperidone) product rows to the medication table,
# right after the existing "APEXIDONE 4mg Tab." Risperidone row (row 144).
# This pushes every subsequent row down by 4 and grows the table from
# A1:J170 to A1:J174.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Grow the table by 4 rows (appended at the bottom by this host). ---
for ($i = 0; $i -lt 4; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# --- 2. Shift existing data rows 145-170 down to 149-174, working from the
#        bottom up so we never overwrite a row before it has been copied. ---
for ($r = 170; $r -ge 145; $r--) {
    $destRow = $r + 4
    $src = $ws.Range("A" + $r + ":J" + $r)
    $dst = $ws.Range("A" + $destRow + ":J" + $destRow)
    $dst.Value = $src.Value()
}

# --- 3. Fill the freed-up rows 145-148 with the new RISPERDAL entries.
#        NOTE: a plain PowerShell @() list does not marshal correctly into
#        a multi-cell Range.Value (COM expects a rectangular object[,]),
#        so each row is built as an explicit 2-D array. ---
function New-RowArray {
    param([object[]] $values)
    $arr = New-Object 'object[,]' 1, $values.Length
    for ($c = 0; $c -lt $values.Length; $c++) {
        $arr[0, $c] = $values[$c]
    }
    return $arr
}

$sideEffects  = "Extrapyramidal symptoms, hyperprolactinemia, weight gain, sedation, orthostatic hypotension"
$interactions = "CNS depressants, drugs that prolong QT interval, levodopa (antagonizes effect), CYP2D6 inhibitors"
$labs         = "Monitor for EPS and akathisia. Dose-dependent EPS risk."
$monitor      = "Sedation, EPS, akathisia, weight gain, prolactin elevation symptoms"
$category     = "Atypical Antipsychotic (Second Generation)"
$halfLife     = "3-20 hours"

$ws.Range("A145:J145").Value = New-RowArray @("RISPERDAL 2mg Tab", "Risperidone", 1, "20TAB", $halfLife, $category, $sideEffects, $interactions, $labs, $monitor)

$ws.Range("A146:J146").Value = New-RowArray @("RISPERDAL 4mg Tab", "Risperidone", 1, "20TAB", $halfLife, $category, $sideEffects, $interactions, $labs, $monitor)

$ws.Range("A147:J147").Value = New-RowArray @("RISPERDAL 1mg/ml Syrup", "Risperidone", 1, "30 ML", $halfLife, $category, $sideEffects, $interactions, $labs, $monitor)

$ws.Range("A148:J148").Value = New-RowArray @("RISPERDAL CONSTA", "Risperidone", 0, "1 SYRING", $halfLife, $category, $sideEffects, $interactions, $labs, $monitor)

# --- 4. Update the sheet view: scroll position + active selection, matching
#        where the newly inserted rows now sit. ---
$sheetView = $ws.Application.ActiveWindow
$ws.Application.ActiveWindow.ScrollRow = 132
$ws.Range("A148").Select()
